$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.847.54"
$ws.Range("E2").Value = "  -0.47%  "

$ws.Range("D3").Value = "1.656.32"
$ws.Range("E3").Value = "  -1.50%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.33"
$ws.Range("E5").Value = "  +0.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3622"
$ws.Range("E7").Value = "  -1.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.13"
$ws.Range("E8").Value = "  -1.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3246"
$ws.Range("E9").Value = "  -3.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.130"
$ws.Range("E10").Value = "  -4.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07046"
$ws.Range("E11").Value = "  -3.95%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.024"
$ws.Range("E13").Value = "  -2.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.46"
$ws.Range("E14").Value = "  -5.32%  "

$ws.Range("D15").Value = "1.657.92"
$ws.Range("E15").Value = "  -1.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.574"
$ws.Range("E16").Value = "  -3.98%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001044"
$ws.Range("E17").Value = "  -5.29%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06576"
$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.78"
$ws.Range("E20").Value = "  -4.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.891"
$ws.Range("E21").Value = "  -5.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.70"
$ws.Range("E22").Value = "  -7.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.52"
$ws.Range("E23").Value = "  -1.09%  "

$ws.Range("D24").Value = "24.836.55"
$ws.Range("E24").Value = "  -0.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.452"
$ws.Range("E26").Value = "  -9.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "147.47"
$ws.Range("E27").Value = "  -2.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.55"
$ws.Range("E28").Value = "  -6.55%  "

$ws.Range("D29").Value = "1.839.07"
$ws.Range("E29").Value = "  -1.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.073"
$ws.Range("E32").Value = "  -2.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.745"
$ws.Range("E33").Value = "  -11.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08445"
$ws.Range("E34").Value = "  -2.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.649"
$ws.Range("E35").Value = "  -4.94%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.24"
$ws.Range("E36").Value = "  -9.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.284"
$ws.Range("E37").Value = "  +2.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.167"
$ws.Range("E38").Value = "  -5.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02254"
$ws.Range("E39").Value = "  -3.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06053"
$ws.Range("E40").Value = "  -6.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.303"
$ws.Range("E41").Value = "  -5.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2065"
$ws.Range("E42").Value = "  -4.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5926"
$ws.Range("E44").Value = "  -5.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.766"
$ws.Range("E45").Value = "  -0.76%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.67"
$ws.Range("E46").Value = "  -6.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5601"
$ws.Range("E47").Value = "  -6.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.03"
$ws.Range("E48").Value = "  -0.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.939"
$ws.Range("E49").Value = "  -5.18%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06969"
$ws.Range("E50").Value = "  -2.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.192"
$ws.Range("E51").Value = "  -0.80%  "

$ws.Range("E25").Value = "  -0.03%  "

# Row 30 and 31: ImmutableX / BitcoinCash swap rank positions
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.196"
$ws.Range("E30").Value = "  -6.75%  "

$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.05"
$ws.Range("E31").Value = "  -4.19%  "
